$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A62").Value = "Casa 🌫️"
$ws.Range("B62").Value = "Christian Martinelli | SdrumALA"
$ws.Range("C62").Value = "Marco Sala | IMONTAGNA"
$ws.Range("D62").Value = "Nadir Chtioui | Mai una gioia"
$ws.Range("E62").Value = "Andrea Roveda | Pinguini Trentini"
$ws.Range("F62").Value = "Stefano Mattioli | SdrumALA"
